$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'66.174.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.75%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.539.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.34%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'580.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.86%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'167.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.37%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.79%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'2.536.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.42%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.05%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  -0.07%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.351"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.37%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'5.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.11%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'26.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.62%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'3.002.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.20%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("E16").Value = "'  -3.17%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'65.991.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.78%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'2.535.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.93%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'11.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -6.83%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.57%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'346.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.77%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'4.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.29%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'4.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.95%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "'  +0.04%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'1.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.03%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'68.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.98%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'9.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -5.25%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  -2.91%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.29%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.0₃0977"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.22%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'525.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.66%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'8.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.91%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  -3.38%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'1.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.92%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.131"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.58%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  +0.27%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("B37").Value = "'ImmutableX"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'1.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.85%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("B38").Value = "'Monero"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'156.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.47%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'18.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.67%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'18.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.68%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.355"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.37%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "'RenderToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'5.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.61%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "'Stacks"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'1.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.93%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("E44").Value = "'  +0.13%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'2.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.88%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "'BabyDogeCoin"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.0₆0282"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.60%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = "'Aave"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'147.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.98%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.556"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.11%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'3.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.42%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'1.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.44%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "'  -1.90%  "
$ws.Range("E51").Style = "Normal"
